$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AH1").Value = "Remarks"
$ws.Range("AH2").Value = "Developer"
$ws.Range("AH3").Value = "Tester"
$ws.Range("AH4").Value = "Tester"

$ws.Range("W2").Value = "Sunil"
$ws.Range("Y2").Value = "Raghu"

$ws.Range("W3").Value = "Pankaj"
$ws.Range("Y3").Value = "Ashok"

$ws.Range("W4").Value = "venkatesh"
$ws.Range("Y4").Value = "Raghu"

$ws.Range("AH4").Select()
